# Auto-generated Word COM-interop script
$d = $word.ActiveDocument

# --- string literals ---
$txt_p1_title = 'San Diego County Gross Regional Product Report'
$txt_p2 = 'Over the past five years, San Diego County has experienced a notable increase in its Gross Regional Product (GRP), reflecting a robust economic growth trajectory. In 2019, the total GRP was $244.28 billion, which saw a slight increase to $244.82 billion in 2020. The following years marked significant growth, with the GRP reaching $268.87 billion in 2021, $296.68 billion in 2022, and ultimately $308.71 billion in 2023. This represents an overall increase of approximately 26.4% from 2019 to 2023, with the most substantial growth occurring between 2021 and 2022, where the GRP rose by about 10.3%.'
$txt_p3 = 'The economic landscape of San Diego County is shaped by several key industries that contribute significantly to its GRP. In 2023, the government sector emerged as the largest contributor with $52.92 billion. This was followed by the professional, scientific, and technical services industry, which added $37.04 billion. Manufacturing also played a crucial role, contributing $31.67 billion, while the health care and social assistance sector added $20.21 billion. The finance and insurance industry rounded out the top contributors with $19.50 billion. These sectors underscore the diverse economic base of San Diego County, highlighting the importance of government, technology, manufacturing, and healthcare in driving regional economic growth.'
$txt_p4 = 'In terms of per capita GRP, San Diego County has consistently outperformed both state and national averages over the past five years. In 2019, the per capita GRP was $73,347, which increased to $74,278 in 2020. The upward trend continued with $82,100 in 2021, $90,557 in 2022, and $94,916 in 2023. When compared to the 2023 figures, San Diego County''s per capita GRP of $94,916 surpasses California''s average of $93,800 and the national average of $77,366. This indicates a relatively prosperous region with a strong economic output per resident, reflecting the high productivity and economic strength of its industries.'
$txt_caption = 'San Diego County GRP and Per Capita GRP (2019-2023)'
$txt_p_concl = 'In conclusion, San Diego County''s economy has demonstrated robust growth over the past five years, driven by significant contributions from the government, professional services, manufacturing, and healthcare sectors. The county''s per capita GRP consistently exceeds both state and national averages, underscoring its economic vitality and the high productivity of its industries. This growth trend positions San Diego County as a key economic hub within California and the United States.'
$txt_footer = 'Page'
$txt_header = 'San Diego County GRP Analysis'

# --- 1. Title paragraph (was Heading1 'Report') ---
$p1 = $d.Paragraphs(1)
$p1.Range.Text = $txt_p1_title
$p1.Style = "Title"

# --- 2. Replace the second paragraph's text ---
$p2 = $d.Paragraphs(2)
$p2.Range.Text = $txt_p2

# --- 3. Insert the 3rd and 4th narrative paragraphs ---
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs(3)
$p3.Range.Text = $txt_p3
$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs(4)
$p4.Range.Text = $txt_p4

# --- 4. Insert the caption paragraph and concluding paragraph after p4 ---
$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs(5)
$p5.Range.Text = $txt_caption
$p5.Style = "Caption"
$p5.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs(6)
$p6.Range.Text = $txt_p_concl

# --- 5. Insert the 6x3 data table immediately before the caption paragraph ---
$tableInsertRange = $p5.Range
$tableInsertRange.Collapse(1)
$table = $d.Tables.Add($tableInsertRange, 6, 3)
$table.Style = "Table Grid"
for ($r=1; $r -le 6; $r++) {
    for ($c=1; $c -le 3; $c++) {
        $table.Cell($r,$c).Range.Style = "Normal"
    }
}

$tableData = @(
    @('Year', 'Total GRP (billion)', 'Per Capita GRP'),
    @('2019', '$244.28', '$73,347'),
    @('2020', '$244.82', '$74,278'),
    @('2021', '$268.87', '$82,100'),
    @('2022', '$296.68', '$90,557'),
    @('2023', '$308.71', '$94,916'),
)
for ($r=1; $r -le 6; $r++) {
    for ($c=1; $c -le 3; $c++) {
        $table.Cell($r,$c).Range.Text = $tableData[$r-1][$c-1]
    }
}

# --- 6. Add the default header and footer ---
$sec = $d.Sections(1)
$hdr = $sec.Headers(1)
$hdr.Range.Text = $txt_header
$ftr = $sec.Footers(1)
$ftr.Range.Text = $txt_footer

Write-Host "Edit complete."
